# "modified all figures, and numbers in text"
#
# The underlying edit: on the "Electricity generation" sheet, the electricity
# price input (row 18, columns C:G) was changed from 966 to 855. Every other
# changed cell in the workbook (rows 25, 27, 31, 37, 40, 44) is a formula that
# depends on row 18, so those values ripple through automatically on recalc.
#
# In addition, the selection on "Electricity generation" moved to C18:G18
# (the cells that were edited) and the active sheet switched from
# "Electricity generation" to "Sheet1".

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Electricity generation")

# Update the input values that drive all the downstream formulas.
$ws2.Range("C18:G18").Value = 855

# Leave the edited range selected on the "Electricity generation" sheet.
$ws2.Range("C18:G18").Select() | Out-Null

# Sheet1 becomes the active/visible sheet.
$ws1.Activate() | Out-Null
